$wb = $excel.ActiveWorkbook

# --- CaliforniaHousing (sheet3) value updates ---
$ws3 = $wb.Worksheets.Item("CaliforniaHousing")
$ws3.Range("L9").Value = 41239.968499999901
$ws3.Range("M9").Value = 28402.999899999999
$ws3.Range("N9").Value = 0.87170000000000003
$ws3.Range("K10").Value = 0.1817
$ws3.Range("L10").Value = 47882.435299999997
$ws3.Range("M10").Value = 32167.810399999998
$ws3.Range("N10").Value = 0.82650000000000001
$ws3.Range("K11").Value = 0.18149999999999999
$ws3.Range("L11").Value = 46073.158600000002
$ws3.Range("M11").Value = 31440.0141999999
$ws3.Range("N11").Value = 0.84569999999999901
$ws3.Range("K14").Value = 0.16170000000000001
$ws3.Range("L14").Value = 41297.653599999998
$ws3.Range("M14").Value = 28367.798299999999
$ws3.Range("N14").Value = 0.87149999999999905
$ws3.Range("K15").Value = 0.18060000000000001
$ws3.Range("L15").Value = 47902.006800000003
$ws3.Range("M15").Value = 32125.017899999901
$ws3.Range("N15").Value = 0.8266
$ws3.Range("K16").Value = 0.17910000000000001
$ws3.Range("L16").Value = 45790.553800000002
$ws3.Range("M16").Value = 31117.5
$ws3.Range("N16").Value = 0.84749999999999903
$ws3.Range("K19").Value = 0.16200000000000001
$ws3.Range("L19").Value = 41390.044199999997
$ws3.Range("M19").Value = 28458.6604999999
$ws3.Range("K20").Value = 0.18149999999999999
$ws3.Range("L20").Value = 48039.911899999999
$ws3.Range("M20").Value = 32261.763299999999
$ws3.Range("N20").Value = 0.82549999999999901
$ws3.Range("K21").Value = 0.18029999999999899
$ws3.Range("L21").Value = 46054.006099999999
$ws3.Range("M21").Value = 31287.5929
$ws3.Range("N21").Value = 0.84570000000000001
$ws3.Range("K24").Value = 0.16189999999999999
$ws3.Range("L24").Value = 41205.103499999997
$ws3.Range("M24").Value = 28393.1611999999
$ws3.Range("N24").Value = 0.87189999999999901
$ws3.Range("K25").Value = 0.183699999999999
$ws3.Range("L25").Value = 48061.705699999999
$ws3.Range("M25").Value = 32379.892699999898
$ws3.Range("N25").Value = 0.82549999999999901
$ws3.Range("K26").Value = 0.18149999999999999
$ws3.Range("L26").Value = 45784.567799999997
$ws3.Range("M26").Value = 31323.0281
$ws3.Range("N26").Value = 0.84749999999999903

# --- Santader (sheet4) value updates ---
$ws4 = $wb.Worksheets.Item("Santader")
$ws4.Range("K14").Value = 4.1393000000000004
$ws4.Range("L14").Value = 3666571.86719999
$ws4.Range("M14").Value = 2350293.1891999999
$ws4.Range("N14").Value = 0.801399999999999
$ws4.Range("K15").Value = 7.5333999999999897
$ws4.Range("L15").Value = 6951613.0497000003
$ws4.Range("M15").Value = 4675228.7275999999
$ws4.Range("N15").Value = 0.28059999999999902
$ws4.Range("K16").Value = 6.5222999999999898
$ws4.Range("L16").Value = 7577015.1265000002
$ws4.Range("M16").Value = 5035427.2412
$ws4.Range("N16").Value = 0.1714
$ws4.Range("K19").Value = 4.4104999999999999
$ws4.Range("L19").Value = 3725080.2409999999
$ws4.Range("M19").Value = 2414190.0791000002
$ws4.Range("N19").Value = 0.79490000000000005
$ws4.Range("K20").Value = 7.7874999999999996
$ws4.Range("L20").Value = 7006900.6945999898
$ws4.Range("M20").Value = 4735679.5045999996
$ws4.Range("N20").Value = 0.26869999999999999
$ws4.Range("K21").Value = 6.8277999999999999
$ws4.Range("L21").Value = 7617075.6844999902
$ws4.Range("M21").Value = 5098734.0138999997
$ws4.Range("N21").Value = 0.16259999999999999
$ws4.Range("K24").Value = 4.4086999999999996
$ws4.Range("L24").Value = 3722339.142
$ws4.Range("M24").Value = 2407725.95299999
$ws4.Range("N24").Value = 0.79500000000000004
$ws4.Range("K25").Value = 7.7796000000000003
$ws4.Range("L25").Value = 6963352.8272000002
$ws4.Range("M25").Value = 4714035.4066000003
$ws4.Range("N25").Value = 0.2777
$ws4.Range("K26").Value = 6.8105999999999902
$ws4.Range("L26").Value = 7619603.5614999998
$ws4.Range("M26").Value = 5109145.7723000003
$ws4.Range("N26").Value = 0.16199999999999901

# --- allstate (sheet5) value updates ---
$ws5 = $wb.Worksheets.Item("allstate")
$ws5.Range("K10").Value = 0.60399999999999898
$ws5.Range("L10").Value = 1771.8401999999901
$ws5.Range("M10").Value = 1159.27529999999
$ws5.Range("N10").Value = 0.62729999999999997
$ws5.Range("K11").Value = 0.61370000000000002
$ws5.Range("L11").Value = 1893.9047
$ws5.Range("M11").Value = 1193.6670999999999
$ws5.Range("N11").Value = 0.57349999999999901
$ws5.Range("K12").Value = 0.60729999999999995
$ws5.Range("L12").Value = 1957.7544
$ws5.Range("M12").Value = 1206.3689999999999
$ws5.Range("K20").Value = 0.60370000000000001
$ws5.Range("L20").Value = 1771.4305999999999
$ws5.Range("M20").Value = 1158.8183999999901
$ws5.Range("N20").Value = 0.62739999999999996
$ws5.Range("K21").Value = 0.61349999999999905
$ws5.Range("L21").Value = 1890.8679
$ws5.Range("M21").Value = 1193.6813
$ws5.Range("N21").Value = 0.57499999999999996
$ws5.Range("K22").Value = 0.60539999999999905
$ws5.Range("L22").Value = 1957.6081999999999
$ws5.Range("M22").Value = 1206.4852000000001
$ws5.Range("N22").Value = 0.55120000000000002

# --- allstate (sheet5): remove obsolete "LOF" duplicate-header cells and the whole
# second (no-imputation-needed) LOF result block that the author dropped ---
$ws5.Range("I9").ClearContents()
$ws5.Range("I19").ClearContents()
$ws5.Range("I24:N27").ClearContents()

# --- sberbank-russian (sheet1): view state (zoom + selection) ---
$ws1 = $wb.Worksheets.Item("sberbank-russian")
$ws1.Activate()
$win = $wb.Windows.Item(1)
$win.Zoom = 55
$ws1.Range("F34").Select()

# --- allstate (sheet5): selection change, and restore as the active/visible tab ---
$ws5.Activate()
$ws5.Range("I37").Select()
